$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-assert a few untouched code samples verbatim so that their
#     embedded carriage returns survive the save round-trip unchanged ---
$ws.Range("C2").Value = "def add(num1, num2):`r`n    return(num1 + num2)"
$ws.Range("C3").Value = "def fibonacci(num):`r`n    if num == 1 or num == 2:`r`n        return num`r`n    else:`r`n        return fibonacci(num - 1) + fibonacci(num - 2)"
$ws.Range("D3").Value = "function fibonacci(num) {`r`n    if (num == 1 || num == 2) {`r`n        return num;`r`n    } else {`r`n        return fibonacci(num-1) + fibonacci(num-2);`r`n    }`r`n}"

# --- Update the code-sample text for the "Hello World" row (row 1) ---
# Python / JavaScript / C# / Java examples now wrap their code in a
# hello_world function (library + main-style function).
$ws.Range("C1").Value = "def hello_world():`n    print('Hello, World!')"
$ws.Range("D1").Value = "function hello_world() {`n     print(`"Hello, World!`");`n}"
$ws.Range("E1").Value = "public static void hello_world()`n{`r`n    Console.WriteLine(`"Hello, World!`");`n}"
$ws.Range("F1").Value = "public static void hello_world() {`n    System.out.print(`"Hello, World!`");`n}"

# --- Update the C# examples for "Adding two numbers" (row 2) and
#     "Fibonacci" (row 3) to use static methods ---
$ws.Range("E2").Value = "public static int add(int num1, int num2)`r`n{`r`n    return num1 + num2;`r`n}"
$ws.Range("F2").Value = "public static int add(int num1, int num2) {`r`n    return num1 + num2;`r`n}"
$ws.Range("E3").Value = "public static int fibonacci(num)`r`n{`r`n    if (num == 1 || num == 2) `r`n    {`r`n        return num;`r`n    }`r`n    else`r`n    {`r`n        return fibonacci(num - 1) + fibonacci(num - 2);`r`n    }`r`n}"
$ws.Range("F3").Value = "public static int fibonacci(int num) {`r`n    if (num == 1 || num == 2) {`r`n        return num;`r`n    } else {`r`n        return fibonacci(num-1) + fibonacci(num-2);`r`n    }`r`n}"

# --- Row 1 now holds multi-line code too, so wrap + size it like the
#     other data rows ---
$ws.Range("C1:F1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 60

# --- Reset the view: no frozen/scrolled left column, and the last
#     selected cell moves from F4 to F3 ---
$ws.Range("F3").Select() | Out-Null
